$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "make monster abilties"
$ws.Range("B3").Value = "not started"
$ws.Range("A4").Value = "make monster attacks"
$ws.Range("B4").Value = "started"
$ws.Range("A5").Value = "make printing functions"
$ws.Range("B5").Value = "not started"
$ws.Range("A6").Value = "make CR calculator"
$ws.Range("B6").Value = "finished?"

$ws.Range("G7").Select()
